# Update quarterly database (simorgh) and refresh quarter-shift values.
#
# The workbook tracks 10 rolling quarters (columns E..N) for a handful of
# expense line items (rows 10-20) and headcount metrics (rows 26-29).
# A new quarter ("فصل اول منتهی به 1401/12") was added to the data set, the
# oldest quarter ("فصل سوم منتهی به 1399/06") was dropped, and every other
# quarter's figures shifted one column to the left (E<-F<-G<-...<-N) with a
# freshly computed value landing in N. The two header rows (8 and 24) that
# label each quarter are updated the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - quarter headers for the expenses table
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل اول منتهی به 1401/12"

# Row 10 - هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 6689
$ws.Range("F10").Value = 6262
$ws.Range("G10").Value = 14279
$ws.Range("H10").Value = 15612
$ws.Range("I10").Value = 11843
$ws.Range("J10").Value = 13049
$ws.Range("K10").Value = 12630
$ws.Range("L10").Value = -25679
$ws.Range("M10").Value = 58084
$ws.Range("N10").Value = 15282

# Row 12 - حق العمل و کمیسیون فروش
$ws.Range("E12").Value = -22171
$ws.Range("F12").Value = 8031
$ws.Range("G12").Value = -8031
$ws.Range("H12").Value = 0

# Row 13 - هزینه تبلیغات
$ws.Range("E13").Value = 838
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 5895
$ws.Range("H13").Value = -5895
$ws.Range("I13").Value = 12277
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 7053
$ws.Range("L13").Value = -7053
$ws.Range("M13").Value = 10630
$ws.Range("N13").Value = 3778

# Row 14 - هزینه مواد مصرفی
$ws.Range("E14").Value = 1273
$ws.Range("F14").Value = 1701
$ws.Range("G14").Value = 1366
$ws.Range("H14").Value = 1798
$ws.Range("I14").Value = 1762
$ws.Range("J14").Value = 1492
$ws.Range("K14").Value = 1915
$ws.Range("L14").Value = -3407
$ws.Range("M14").Value = 7923
$ws.Range("N14").Value = 2642

# Row 15 - هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 1574
$ws.Range("F15").Value = 1450
$ws.Range("G15").Value = 1416
$ws.Range("H15").Value = 1777
$ws.Range("I15").Value = 1723
$ws.Range("J15").Value = 2021
$ws.Range("K15").Value = 1121
$ws.Range("L15").Value = -3142
$ws.Range("M15").Value = 7877
$ws.Range("N15").Value = 2597

# Row 16 - هزینه استهلاک
$ws.Range("E16").Value = 3123
$ws.Range("F16").Value = 2999
$ws.Range("G16").Value = 2643
$ws.Range("H16").Value = 2941
$ws.Range("I16").Value = 5597
$ws.Range("J16").Value = 39201
$ws.Range("K16").Value = -30448
$ws.Range("L16").Value = -8753
$ws.Range("M16").Value = 18690
$ws.Range("N16").Value = 6693

# Row 17 - هزینه حقوق و دستمزد
$ws.Range("E17").Value = 81943
$ws.Range("F17").Value = 71309
$ws.Range("G17").Value = 77261
$ws.Range("H17").Value = 78509
$ws.Range("I17").Value = 269853
$ws.Range("J17").Value = 110977
$ws.Range("K17").Value = 132394
$ws.Range("L17").Value = -243371
$ws.Range("M17").Value = 763265
$ws.Range("N17").Value = 197376

# Row 19 - سایر هزینه ها
$ws.Range("E19").Value = 147053
$ws.Range("F19").Value = 58608
$ws.Range("G19").Value = 80974
$ws.Range("H19").Value = 93127
$ws.Range("I19").Value = -40090
$ws.Range("J19").Value = 66595
$ws.Range("K19").Value = 197070
$ws.Range("L19").Value = -263665
$ws.Range("M19").Value = 395911
$ws.Range("N19").Value = 97948

# Row 20 - جمع
$ws.Range("E20").Value = 220322
$ws.Range("F20").Value = 150360
$ws.Range("G20").Value = 175803
$ws.Range("H20").Value = 187869
$ws.Range("I20").Value = 262965
$ws.Range("J20").Value = 233335
$ws.Range("K20").Value = 321735
$ws.Range("L20").Value = -555070
$ws.Range("M20").Value = 1262380
$ws.Range("N20").Value = 326316

# Row 24 - quarter headers for the headcount table
$ws.Range("E24").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F24").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G24").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H24").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I24").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J24").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K24").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L24").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M24").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N24").Value = "فصل اول منتهی به 1401/12"

# Row 26 - تعداد پرسنل شرکت در بخش های اداری و پشتیبانی
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = "-"
$ws.Range("K26").Value = 128
$ws.Range("N26").Value = 128

# Row 27 - تعداد پرسنل شرکت در بخش تولید
$ws.Range("I27").Value = 1571
$ws.Range("J27").Value = "-"
$ws.Range("K27").Value = 1560
$ws.Range("M27").Value = 1502
$ws.Range("N27").Value = 1502

# Row 28 - تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E28").Value = 47
$ws.Range("F28").Value = 45
$ws.Range("G28").Value = 44
$ws.Range("H28").Value = 47
$ws.Range("I28").Value = "-"
$ws.Range("J28").Value = 43
$ws.Range("K28").Value = "-"

# Row 29 - تعداد پرسنل تولیدی شرکت
$ws.Range("E29").Value = 1607
$ws.Range("F29").Value = 1685
$ws.Range("G29").Value = 1664
$ws.Range("H29").Value = 1607
$ws.Range("I29").Value = "-"
$ws.Range("J29").Value = 1656
$ws.Range("K29").Value = "-"
